$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 555, shifting existing rows 555-606 down to 556-607
$ws.Rows.Item(555).Insert()

# Populate the newly inserted row 555 with the new weekly data entry
$row = 555
$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 45223
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112032
$ws.Cells.Item($row, 7).Value = "Zapallo italiano"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 440
$ws.Cells.Item($row, 11).Value = 11500
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 11750
$ws.Cells.Item($row, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 196
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
